$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 / Row 11 swap (Dogecoin <-> Toncoin) with updated values ---
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.45"
$ws.Range("E10").Value = "  -3.28%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  +3.63%  "

# --- Row 40 / Row 41 swap (TheGraph <-> OKB) with updated values ---
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'50.14"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  +3.46%  "

# --- Price (column D) updates ---
# The apostrophe prefix forces Excel to keep these values as text (matching the
# original inlineStr cell type) instead of auto-converting them to numbers.
$ws.Range("D2").Value = "'66.863.75"
$ws.Range("D3").Value = "'3.099.55"
$ws.Range("D5").Value = "'580.79"
$ws.Range("D6").Value = "'172.71"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D8").Value = "'3.095.03"
$ws.Range("D12").Value = "'0.484"
$ws.Range("D14").Value = "'37.45"
$ws.Range("D16").Value = "'3.614.44"
$ws.Range("D17").Value = "'66.866.34"
$ws.Range("D19").Value = "'3.102.19"
$ws.Range("D20").Value = "'16.24"
$ws.Range("D21").Value = "'482.27"
$ws.Range("D24").Value = "'84.01"
$ws.Range("D25").Value = "'2.37"
$ws.Range("D26").Value = "'13.12"
$ws.Range("D27").Value = "'10.05"
$ws.Range("D29").Value = "'7.97"
$ws.Range("D33").Value = "'28.65"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D36").Value = "'5.90"
$ws.Range("D37").Value = "'0.991"
$ws.Range("D38").Value = "'48.05"
$ws.Range("D39").Value = "'2.11"
$ws.Range("D44").Value = "'2.81"
$ws.Range("D46").Value = "'2.836.37"
$ws.Range("D47").Value = "'384.33"
$ws.Range("D48").Value = "'134.61"
$ws.Range("D50").Value = "'24.95"
$ws.Range("D51").Value = "'2.23"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("E3").Value = "  +5.08%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("E6").Value = "  +6.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +5.05%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("E14").Value = "  +7.49%  "
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +5.09%  "
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("E21").Value = "  +8.05%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("E25").Value = "  +5.26%  "
$ws.Range("E26").Value = "  +6.90%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("E33").Value = "  +5.23%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("E46").Value = "  +5.84%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("E51").Value = "  +2.68%  "
